$wb = $excel.ActiveWorkbook

# --- Sheet1 "Generic": NrBuckets 4 -> 5 ---
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 5

# --- Sheet3 "Productdata" ---
$wsProd = $wb.Worksheets.Item("Productdata")
$wsProd.Range("E2").Value = 11.6928064
$wsProd.Range("E3").Value = 2.1893696
$wsProd.Range("E4").Value = 0.8987520000000001
$wsProd.Range("E5").Value = 1.447344
$wsProd.Range("E6").Value = 0.9553424
$wsProd.Range("E7").Value = 0.292656
$wsProd.Range("E8").Value = 0.09440640000000002
$wsProd.Range("E9").Value = 0.8622864
$wsProd.Range("E10").Value = 0.5180544
$wsProd.Range("E11").Value = 0.757008
$wsProd.Range("C12").Value = 465
$wsProd.Range("E12").Value = 1.2475936
$wsProd.Range("C13").Value = 3745
$wsProd.Range("E13").Value = 12.57430399999999
$wsProd.Range("C14").Value = 1394
$wsProd.Range("E14").Value = 4.885091200000001
$wsProd.Range("C15").Value = 342
$wsProd.Range("E15").Value = 0.8773248000000001
$wsProd.Range("C16").Value = 487
$wsProd.Range("E16").Value = 0.8904128
$wsProd.Range("C17").Value = 684
$wsProd.Range("E17").Value = 1.393704
$wsProd.Range("C18").Value = 204
$wsProd.Range("E18").Value = 0.45864
$wsProd.Range("C19").Value = 57
$wsProd.Range("E19").Value = 0.1306432
$wsProd.Range("E20").Value = 63.16628800000002
$wsProd.Range("E21").Value = 67.05713920000001
$wsProd.Range("E22").Value = 83.00482560000002
$wsProd.Range("E23").Value = 255.4031584000001

# --- Sheet4 "ForecastedAverageDemand" ---
$wsFAD = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsFAD.Range("C2").Value = 601
$wsFAD.Range("D2").Value = 148
$wsFAD.Range("F2").Value = 302
$wsFAD.Range("G2").Value = 91
$wsFAD.Range("H2").Value = 28
$wsFAD.Range("I2").Value = 418
$wsFAD.Range("J2").Value = 296
$wsFAD.Range("K2").Value = 906
$wsFAD.Range("L2").Value = 210
$wsFAD.Range("C3").Value = 607
$wsFAD.Range("F3").Value = 292
$wsFAD.Range("G3").Value = 94
$wsFAD.Range("H3").Value = 34
$wsFAD.Range("I3").Value = 422
$wsFAD.Range("J3").Value = 300
$wsFAD.Range("K3").Value = 898
$wsFAD.Range("L3").Value = 211
$wsFAD.Range("C4").Value = 610
$wsFAD.Range("D4").Value = 154
$wsFAD.Range("F4").Value = 295
$wsFAD.Range("G4").Value = 91
$wsFAD.Range("I4").Value = 420
$wsFAD.Range("K4").Value = 907
$wsFAD.Range("L4").Value = 213
$wsFAD.Range("C5").Value = 600
$wsFAD.Range("D5").Value = 143
$wsFAD.Range("F5").Value = 299
$wsFAD.Range("G5").Value = 86
$wsFAD.Range("H5").Value = 21
$wsFAD.Range("I5").Value = 418
$wsFAD.Range("J5").Value = 302
$wsFAD.Range("K5").Value = 900
$wsFAD.Range("L5").Value = 210
# new row 6 (copy style of row 5 column A first)
$wsFAD.Range("A5").Copy($wsFAD.Range("A6"))
$wsFAD.Range("A6").Value = 4
$wsFAD.Range("B6").Value = 0
$wsFAD.Range("C6").Value = 596
$wsFAD.Range("D6").Value = 150
$wsFAD.Range("E6").Value = 0
$wsFAD.Range("F6").Value = 301
$wsFAD.Range("G6").Value = 93
$wsFAD.Range("H6").Value = 37
$wsFAD.Range("I6").Value = 419
$wsFAD.Range("J6").Value = 299
$wsFAD.Range("K6").Value = 895
$wsFAD.Range("L6").Value = 214
$wsFAD.Range("M6").Value = 0
$wsFAD.Range("N6").Value = 0
$wsFAD.Range("O6").Value = 0
$wsFAD.Range("P6").Value = 0
$wsFAD.Range("Q6").Value = 0
$wsFAD.Range("R6").Value = 0
$wsFAD.Range("S6").Value = 0
$wsFAD.Range("T6").Value = 0
$wsFAD.Range("U6").Value = 0
$wsFAD.Range("V6").Value = 0
$wsFAD.Range("W6").Value = 0

# --- Sheet5 "ForcastedStandardDeviation" ---
$wsFSD = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsFSD.Range("C2").Value = 75.125
$wsFSD.Range("D2").Value = 18.5
$wsFSD.Range("F2").Value = 37.75
$wsFSD.Range("G2").Value = 11.375
$wsFSD.Range("H2").Value = 3.5
$wsFSD.Range("I2").Value = 52.25
$wsFSD.Range("J2").Value = 37
$wsFSD.Range("K2").Value = 113.25
$wsFSD.Range("L2").Value = 26.25
$wsFSD.Range("C3").Value = 113.8125
$wsFSD.Range("F3").Value = 54.75
$wsFSD.Range("G3").Value = 17.625
$wsFSD.Range("H3").Value = 6.375
$wsFSD.Range("I3").Value = 79.125
$wsFSD.Range("J3").Value = 56.25
$wsFSD.Range("K3").Value = 168.375
$wsFSD.Range("L3").Value = 39.5625
$wsFSD.Range("C4").Value = 133.4375
$wsFSD.Range("D4").Value = 33.6875
$wsFSD.Range("F4").Value = 64.53125
$wsFSD.Range("G4").Value = 19.90625
$wsFSD.Range("I4").Value = 91.875
$wsFSD.Range("K4").Value = 198.40625
$wsFSD.Range("L4").Value = 46.59375
$wsFSD.Range("C5").Value = 140.625
$wsFSD.Range("D5").Value = 33.515625
$wsFSD.Range("F5").Value = 70.078125
$wsFSD.Range("G5").Value = 20.15625
$wsFSD.Range("H5").Value = 4.921875
$wsFSD.Range("I5").Value = 97.96875
$wsFSD.Range("J5").Value = 70.78125
$wsFSD.Range("K5").Value = 210.9375
$wsFSD.Range("L5").Value = 49.21875
# new row 6 (copy style of row 5 column A first)
$wsFSD.Range("A5").Copy($wsFSD.Range("A6"))
$wsFSD.Range("A6").Value = 4
$wsFSD.Range("B6").Value = 0
$wsFSD.Range("C6").Value = 144.34375
$wsFSD.Range("D6").Value = 36.328125
$wsFSD.Range("E6").Value = 0
$wsFSD.Range("F6").Value = 72.8984375
$wsFSD.Range("G6").Value = 22.5234375
$wsFSD.Range("H6").Value = 8.9609375
$wsFSD.Range("I6").Value = 101.4765625
$wsFSD.Range("J6").Value = 72.4140625
$wsFSD.Range("K6").Value = 216.7578125
$wsFSD.Range("L6").Value = 51.828125
$wsFSD.Range("M6").Value = 0
$wsFSD.Range("N6").Value = 0
$wsFSD.Range("O6").Value = 0
$wsFSD.Range("P6").Value = 0
$wsFSD.Range("Q6").Value = 0
$wsFSD.Range("R6").Value = 0
$wsFSD.Range("S6").Value = 0
$wsFSD.Range("T6").Value = 0
$wsFSD.Range("U6").Value = 0
$wsFSD.Range("V6").Value = 0
$wsFSD.Range("W6").Value = 0

# --- Sheet6 "Capacity" ---
$wsCap = $wb.Worksheets.Item("Capacity")
$wsCap.Range("B2").Value = 50344
$wsCap.Range("B3").Value = 525385
$wsCap.Range("B4").Value = 525385
$wsCap.Range("B5").Value = 4953630
